$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched columns remain plain text so values like "27.126.98"
# or leading/trailing-space percentages are not reinterpreted by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$updates = @{
    "D2"  = "27.126.98"
    "E2"  = "  +2.41%  "
    "D3"  = "1.780.13"
    "E3"  = "  +3.57%  "
    "D4"  = "1.015"
    "E4"  = "  +1.18%  "
    "D5"  = "334.42"
    "E5"  = "  +0.49%  "
    "D6"  = "1.009"
    "E6"  = "  +0.95%  "
    "D7"  = "0.3765"
    "E7"  = "  +1.81%  "
    "B8"  = "Cardano"
    "C8"  = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
    "D8"  = "0.3417"
    "E8"  = "  +2.07%  "
    "B9"  = "OKB"
    "C9"  = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D9"  = "48.45"
    "E9"  = "  +0.61%  "
    "D10" = "1.193"
    "E10" = "  +0.95%  "
    "D11" = "0.07433"
    "E11" = "  +0.85%  "
    "D12" = "1.013"
    "E12" = "  +1.09%  "
    "D13" = "21.67"
    "E13" = "  +8.19%  "
    "D14" = "6.438"
    "E14" = "  +0.92%  "
    "D15" = "1.788.14"
    "E15" = "  +4.10%  "
    "D16" = "7.001"
    "E16" = "  -0.28%  "
    "D17" = "0.00001086"
    "E17" = "  +1.84%  "
    "D18" = "0.06644"
    "E18" = "  +0.38%  "
    "D19" = "83.91"
    "E19" = "  +2.46%  "
    "D20" = "1.006"
    "E20" = "  +0.57%  "
    "D21" = "17.21"
    "E21" = "  +4.39%  "
    "D22" = "6.407"
    "E22" = "  +4.70%  "
    "D23" = "27.164.20"
    "E23" = "  +2.60%  "
    "D24" = "12.34"
    "E24" = "  -3.18%  "
    "D25" = "2.460"
    "E25" = "  +1.14%  "
    "D26" = "2.520"
    "E26" = "  +5.80%  "
    "D27" = "1.470"
    "E27" = "  +5.74%  "
    "D28" = "21.28"
    "E28" = "  +10.00%  "
    "D29" = "149.18"
    "E29" = "  -1.84%  "
    "D30" = "1.990.83"
    "E30" = "  +4.27%  "
    "D31" = "132.74"
    "E31" = "  +1.60%  "
    "D32" = "4.067"
    "E32" = "  -1.09%  "
    "D33" = "5.994"
    "E33" = "  +1.53%  "
    "D34" = "0.08627"
    "E34" = "  +0.37%  "
    "D35" = "12.98"
    "E35" = "  +3.08%  "
    "D36" = "1.661"
    "E36" = "  -2.31%  "
    "D37" = "5.388"
    "E37" = "  +1.23%  "
    "D38" = "0.6818"
    "E38" = "  +10.88%  "
    "B39" = "Hedera"
    "C39" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D39" = "0.06290"
    "E39" = "  +2.00%  "
    "B40" = "Algorand"
    "C40" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D40" = "0.2182"
    "E40" = "  +1.43%  "
    "B41" = "VeChain"
    "C41" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D41" = "0.02326"
    "E41" = "  +0.42%  "
    "D42" = "8.732"
    "E42" = "  +3.82%  "
    "D43" = "1.261"
    "E43" = "  +3.31%  "
    "D44" = "14.29"
    "E44" = "  +1.08%  "
    "B45" = "Frax"
    "C45" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
    "D45" = "1.007"
    "E45" = "  +0.80%  "
    "B46" = "Decentraland"
    "C46" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D46" = "0.6372"
    "E46" = "  +7.12%  "
    "D47" = "3.838"
    "E47" = "  -1.63%  "
    "D48" = "2.105"
    "E48" = "  +3.50%  "
    "D49" = "128.93"
    "E49" = "  +0.89%  "
    "D50" = "0.07178"
    "E50" = "  +0.39%  "
    "D51" = "78.89"
    "E51" = "  +2.90%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
